$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.375.02"
$ws.Range("E2").Value = "  +0.28%  "

# Row 3
$ws.Range("D3").Value = "1.876.89"
$ws.Range("E3").Value = "  +0.63%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7124"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.76%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.45"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.59%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07825"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.08%  "

# Row 9
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3116"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.50%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.05"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.61%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08478"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.84%  "

# Row 12
$ws.Range("D12").Value = "1.882.54"
$ws.Range("E12").Value = "  -8.53%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.242"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.14%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7132"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.00%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.56"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.73%  "

# Row 16
$ws.Range("D16").Value = "29.363.06"
$ws.Range("E16").Value = "  +0.12%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008253"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.60%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.054"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.98%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.18"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.27%  "

# Row 20
$ws.Range("E20").Value = "  +0.76%  "

# Row 21
$ws.Range("D21").Value = "2.127.71"
$ws.Range("E21").Value = "  -0.19%  "

# Row 22
$ws.Range("E22").Value = "  -0.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.836"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.51%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.06%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1594"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.45%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.62"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.065"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.31%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.52"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.98%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.516"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.29%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.435"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.44%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.349"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.52%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.282"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.69%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05336"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.49%  "

# Row 34
$ws.Range("E34").Value = "  +0.33%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7562"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.91%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.179"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.28%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.687"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.02%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01875"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.62%  "

# Row 39
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.210.13"
$ws.Range("E39").Value = "  +2.05%  "

# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.722"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.78%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.486"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.63%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.18"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.09%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8869"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.31%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "107.98"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.76%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.02%  "

# Row 46
$ws.Range("D46").Value = "2.025.16"
$ws.Range("E46").Value = "  +0.20%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.825"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.48%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5209"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.24%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000123"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +10.86%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.451"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.74%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4331"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.14%  "

Write-Output "Applied crypto price updates"